$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins, Losses, Ties in AD1:AF1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold, border, centered) from A1 onto
# the three new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-39: every row gets the same team record values (93 wins,
# 69 losses, 0 ties).
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 30).Value = 93   # AD
    $ws.Cells.Item($r, 31).Value = 69   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
